$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 15.346
$ws.Range("C2").Value = 8.702
$ws.Range("D2").Value = 24.048

$ws.Range("B3").Value = 23.623
$ws.Range("C3").Value = 47.127
$ws.Range("D3").Value = 70.75

$ws.Range("B4").Value = 1.726
$ws.Range("C4").Value = 0.875
$ws.Range("D4").Value = 2.601

$ws.Range("C5").Value = 2.601
$ws.Range("D5").Value = 2.601

$ws.Range("B6").Value = 40.695
$ws.Range("C6").Value = 59.305
